$d = $word.ActiveDocument

# Word-style Find/Replace wildcard matching is off; we do plain literal matches.
# MatchCase=$true to avoid any case folding surprises (irrelevant for CJK, harmless).

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# ---------------------------------------------------------------------------
# Paragraph 1: expand the sentence about the book's structure.
# ---------------------------------------------------------------------------

# Insert the parenthetical list of thinkers right after "传记"
Replace-Text "传记，他们" "传记（彼特拉克、马丁·路德、笛卡尔、霍布斯、康德），他们"

# Split "...关键人物，但作者..." into two sentences.
Replace-Text "关键人物，但作者" "关键人物。但作者"

# "再写哲学史" -> "在写哲学史"
Replace-Text "写作方式却像是再写哲学史" "写作方式却像是在写哲学史"

# "思想者本人地生活" -> "思想者本人的生活"
Replace-Text "思想者本人地生活" "思想者本人的生活"

# Insert "方面" before "展开探讨"
Replace-Text "活动等展开探讨" "活动等方面展开探讨"

# "母题" -> "主题"
Replace-Text "这一母题串联起来" "这一主题串联起来"

# "并得出作者的结论" -> "并得出最终的结论"
Replace-Text "并得出作者的结论" "并得出最终的结论"

# ---------------------------------------------------------------------------
# Paragraph 2: the author's definition of "modernity".
# ---------------------------------------------------------------------------

# "并不肯定、明确" -> "并不清晰、明确"
Replace-Text "并不肯定、明确" "并不清晰、明确"

# Move closing quote before the period, add a dun-comma between the two quotes
Replace-Text "相互作用的产物。”“现代性" "相互作用的产物”、“现代性"

# Move the closing quote before the trailing period of the second quotation
Replace-Text "前所未有的。”作者对" "前所未有的”。作者对"

# Add the new clause about a "self-realization" starting point
Replace-Text "有利的限制。作者最后" "有利的限制，一个“自我实现”的出发点。作者最后"

# ---------------------------------------------------------------------------
# Paragraph 4: knowledge background / scope of the critique.
# ---------------------------------------------------------------------------

# Insert parenthetical aside about the author's lack of awareness
Replace-Text "的知识背景造成的" "的知识背景（至少就行文中，未能看到作者对此的认识）造成的"

# "包括科学在内的与" -> "包括科学等与"
Replace-Text "包括科学在内的与" "包括科学等与"

# "都未得到考量" -> "都未纳入作者的考量"
Replace-Text "紧密相连的问题都未得到考量" "紧密相连的问题都未纳入作者的考量"

# ---------------------------------------------------------------------------
# Note: the source diff also drops the (purely cosmetic, non-visible)
# w:hint="eastAsia" attribute from the 4th paragraph's paragraph-mark run
# properties (w:pPr/w:rPr/w:rFonts). That attribute is an internal rendering
# hint that is not exposed anywhere on the Word object model (Range.Font,
# ParagraphFormat, etc. have no such member), so it cannot be toggled via
# COM automation here; it has no effect on the document's visible content
# or formatting.
# ---------------------------------------------------------------------------
